$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 41 (row 50): KKR vs DC
$ws.Range("E50").Value = 60
$ws.Range("H50").Value = 40
$ws.Range("K50").Value = 80
$ws.Range("N50").Value = 20
$ws.Range("Q50").Value = 0
$ws.Range("T50").Value = 100

# Contest 42 (row 51): MI vs PBKS
$ws.Range("E51").Value = 80
$ws.Range("H51").Value = 20
$ws.Range("K51").Value = 60
$ws.Range("N51").Value = 100
$ws.Range("Q51").Value = 0
$ws.Range("T51").Value = 40
